$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daten")
if (-not $ws) { $ws = $wb.ActiveSheet }

# New column I header: "commercial" (constraint 7 - commercial buildings min. distance)
$ws.Range("I1").Value = "commercial"

# Row 2 (Hotel) updates
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 0.01
$ws.Range("E2").Value = 50000
$ws.Range("H2").Value = 0.2
$ws.Range("I2").Value = "'True"

# Row 3 (Theatre) - new column I value
$ws.Range("I3").Value = "'False"

# Row 4 (Movie Theater) - new column I value
$ws.Range("I4").Value = "'False"

# Row 5 (Hospital) updates
$ws.Range("B5").Value = 13000
$ws.Range("I5").Value = "'False"

# Make the new True/False text cells pick up the same (column default) style
# as the header cell instead of the quote-prefixed style COM applies by default.
$ws.Range("I1").Copy()
$ws.Range("I2:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect final cursor position
$ws.Range("E15").Select()
